# Fruta / hortaliza, semanal
# Rotates the weekly price-record data (columns D, L-T) among rows 2,3,4,5,6,7,8,10
# Row 9 is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, taken from the target (post-edit) state.
$rows = @{
    2  = @{ D = 44495; L = "Primera";  M = 50;  N = 24000; O = 24000; P = 24000; Q = "`$/bandeja 10 kilos";       R = "China";                  S = 2400; T = 10 }
    3  = @{ D = 44418; L = "Especial"; M = 100; N = 8000;  O = 8000;  P = 8000;  Q = "`$/caja 15 kilos granel";    R = "Región de O'Higgins";    S = 533;  T = 15 }
    4  = @{ D = 44208; L = "Especial"; M = 70;  N = 24000; O = 24000; P = 24000; Q = "`$/caja 15 kilos granel";    R = "Región de O'Higgins";    S = 1600; T = 15 }
    5  = @{ D = 44427; L = "Primera";  M = 55;  N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 15 kilos granel";    R = "Región de O'Higgins";    S = 467;  T = 15 }
    6  = @{ D = 44411; L = "Primera";  M = 210; N = 8000;  O = 8000;  P = 8000;  Q = "`$/bandeja 8 kilos";         R = "Región de O'Higgins";    S = 1000; T = 8  }
    7  = @{ D = 44217; L = "Primera";  M = 55;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos granel";    R = "Región de O'Higgins";    S = 1000; T = 18 }
    8  = @{ D = 44511; L = "Primera";  M = 15;  N = 22000; O = 22000; P = 22000; Q = "`$/caja 15 kilos granel";    R = "Región de O'Higgins";    S = 1467; T = 15 }
    10 = @{ D = 44392; L = "Especial"; M = 500; N = 7000;  O = 7000;  P = 7000;  Q = "`$/bandeja 8 kilos";         R = "Región de O'Higgins";    S = 875;  T = 8  }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
